# Big-Data-Integration-Paper.docx update ("update to the paper")
#
# 1. Author byline: add four more co-authors after "Sidney Schaeper".
# 2. Under the "Understanding Our Topic" heading, add an intro paragraph,
#    a 3-item bullet list of research questions, and a closing paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Extend the author line:
#    "Sidney Schaeper" -> "Sidney Schaeper, Aman Rastogi, Keya Satpathy,
#    Bhawna Saini, Prathiba Swamykannu"
# ---------------------------------------------------------------------
$authorPara = $d.Paragraphs(2)

$authorAdditions = @(
  ",",
  " ",
  "Aman",
  " ",
  "Rastogi,",
  " ",
  "Keya",
  " ",
  "Satpathy,",
  " ",
  "Bhawna",
  " ",
  "Saini,",
  " ",
  "Prathiba",
  " ",
  "Swamykannu"
)

foreach ($piece in $authorAdditions) {
  $authorPara.Range.InsertAfter($piece)
}

# ---------------------------------------------------------------------
# 2. Insert the new "research questions" content right after the
#    "Understanding Our Topic" heading paragraph.
# ---------------------------------------------------------------------
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $candidate = $d.Paragraphs($i)
  if ($candidate.Range.Text.TrimEnd([char]13, [char]7) -eq "Understanding Our Topic") {
    $headingPara = $candidate
    break
  }
}

$headingPara.Range.InsertParagraphAfter()
$introPara = $headingPara.Next()
$introPara.Style = "FirstParagraph"
$introPara.Range.Text = "We had a variety of questions we wanted to answer in this project. We formualted these questions after evalauting our 84.51 datasets and the twitter data. The following are the questions we hoped to answer."

$introPara.Range.InsertParagraphAfter()
$bullet1 = $introPara.Next()
$bullet1.Style = "Compact"
$bullet1.Range.Text = "Should we continue to do business with our five worst performing brands based on the perspective of customers?"

$bullet1.Range.InsertParagraphAfter()
$bullet2 = $bullet1.Next()
$bullet2.Style = "Compact"
$bullet2.Range.Text = "Does location of the product within the store and ad impact sales of the product? What is the optimal location within the store and ads?"

$bullet2.Range.InsertParagraphAfter()
$bullet3 = $bullet2.Next()
$bullet3.Style = "Compact"
$bullet3.Range.Text = "What amounts of each product should we plan to have in inventory on average daily based on these sales?"

# Apply bullet-list numbering to all three questions in one shot so they
# share a single abstractNum/num definition (numId 1001).
$listRange = $d.Range($bullet1.Range.Start, $bullet3.Range.End)
$listRange.ListFormat.ApplyBulletDefault()

$bullet3.Range.InsertParagraphAfter()
$closingPara = $bullet3.Next()
$closingPara.Style = "FirstParagraph"
$closingPara.Range.Text = "For the first question, we thought that the twitter dataset, product_lookup dataset, and transactions dataset would provide us insight into answering this question. For the second question, we thought that the transactions dataset and casual_lookup dataset would provide us information into answering this question. For the third question, we thought that the transactions dataset, store_lookup dataset, and product_lookup dataset would help us discover the answer to this question."

Write-Output "edit complete"
